$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph (currently the 2nd paragraph,
#    right after the Heading1 title paragraph).
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. Replace the "Create a cartoon-style feature image..." paragraph text with the
#    new meta-description text, keeping its existing italic run formatting.
$d.Content.Find.Execute(
    "Create a cartoon-style feature image for " + [char]34 + "Gods of Giza" + [char]34 + " that features a happy Maya warrior. The image should be eye-catching and fun, with the warrior wearing glasses and smiling to invite players to try their luck on this exciting game. The background should incorporate elements of Ancient Egypt, such as pyramids and hieroglyphs, to give players a glimpse of the game's unique theme. Use bright colors and bold lines to make the image pop and stand out from other slot games. The image should capture the adventurous spirit of the game and convey a sense of excitement and anticipation for players.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Want to play Gods of Giza for free? Read our review of this enhanced slot machine featuring rotating reels, scatters, and free spins.",
    2
) | Out-Null

# 3. Insert a new bold paragraph right before that paragraph with the title text.
$target = $d.Content.Find
$titleRange = $d.Content
$titleRange.Find.Execute(
    "Want to play Gods of Giza for free? Read our review of this enhanced slot machine featuring rotating reels, scatters, and free spins.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null

$para = $titleRange.Paragraphs(1)
$newPara = $para.Range.InsertParagraphBefore()

$insertedRange = $para.Range
